# Generate Report for Handoff
# - Flip the localization status from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it appears (Overview + per-language sheets).
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps to reflect the new handoff generation run.
# - The Status column got a lot narrower once the text shrank, so shrink the
#   column to match (closest value the ColumnWidth pixel grid can represent).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newWidth  = 16.333333333333332   # -> stored col width ~17.1667 (closest to 17.2159881591797)

# ---- Overview sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = $newStatus
$ws.Range("F2").Value = $newStatus
$ws.Range("G2").Value = "2016-08-13 03:12:47"
$ws.Columns.Item(5).ColumnWidth = $newWidth
$ws.Columns.Item(6).ColumnWidth = $newWidth

# ---- zh-cn sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = "2016-08-13 03:12:40"
$ws.Columns.Item(3).ColumnWidth = $newWidth

# ---- de-de sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = $newStatus
$ws.Range("H2").Value = "2016-08-13 03:12:47"
$ws.Columns.Item(3).ColumnWidth = $newWidth
